$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 1144, shifting the rest of the
# table (old rows 1144-1176) down to rows 1147-1179.
$ws.Rows("1144:1146").Insert()

# ---- New row 1144: Palta Hass "Especial" @ Provincia de Limari ----
$ws.Cells.Item(1144,1).Value  = 8
$ws.Cells.Item(1144,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1144,3).Value  = "Coquimbo"
$ws.Cells.Item(1144,4).Value  = 44939
$ws.Cells.Item(1144,5).Value  = 4
$ws.Cells.Item(1144,6).Value  = "Fruta"
$ws.Cells.Item(1144,7).Value  = 100106
$ws.Cells.Item(1144,8).Value  = "Oleaginosos"
$ws.Cells.Item(1144,9).Value  = 100106002
$ws.Cells.Item(1144,10).Value = "Palta"
$ws.Cells.Item(1144,11).Value = "Hass"
$ws.Cells.Item(1144,12).Value = "Especial"
$ws.Cells.Item(1144,13).Value = 160
$ws.Cells.Item(1144,14).Value = 3200
$ws.Cells.Item(1144,15).Value = 3300
$ws.Cells.Item(1144,16).Value = 3250
$ws.Cells.Item(1144,17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(1144,18).Value = "Provincia de Limar$([char]0xED)"
$ws.Cells.Item(1144,19).Value = 3250
$ws.Cells.Item(1144,20).Value = 1

# ---- New row 1145: Palta Hass "Primera" @ Provincia de Limari ----
$ws.Cells.Item(1145,1).Value  = 8
$ws.Cells.Item(1145,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1145,3).Value  = "Coquimbo"
$ws.Cells.Item(1145,4).Value  = 44939
$ws.Cells.Item(1145,5).Value  = 4
$ws.Cells.Item(1145,6).Value  = "Fruta"
$ws.Cells.Item(1145,7).Value  = 100106
$ws.Cells.Item(1145,8).Value  = "Oleaginosos"
$ws.Cells.Item(1145,9).Value  = 100106002
$ws.Cells.Item(1145,10).Value = "Palta"
$ws.Cells.Item(1145,11).Value = "Hass"
$ws.Cells.Item(1145,12).Value = "Primera"
$ws.Cells.Item(1145,13).Value = 200
$ws.Cells.Item(1145,14).Value = 2900
$ws.Cells.Item(1145,15).Value = 3000
$ws.Cells.Item(1145,16).Value = 2950
$ws.Cells.Item(1145,17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(1145,18).Value = "Provincia de Limar$([char]0xED)"
$ws.Cells.Item(1145,19).Value = 2950
$ws.Cells.Item(1145,20).Value = 1

# ---- New row 1146: Palta Hass "Segunda" @ Provincia de Limari ----
$ws.Cells.Item(1146,1).Value  = 8
$ws.Cells.Item(1146,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1146,3).Value  = "Coquimbo"
$ws.Cells.Item(1146,4).Value  = 44939
$ws.Cells.Item(1146,5).Value  = 4
$ws.Cells.Item(1146,6).Value  = "Fruta"
$ws.Cells.Item(1146,7).Value  = 100106
$ws.Cells.Item(1146,8).Value  = "Oleaginosos"
$ws.Cells.Item(1146,9).Value  = 100106002
$ws.Cells.Item(1146,10).Value = "Palta"
$ws.Cells.Item(1146,11).Value = "Hass"
$ws.Cells.Item(1146,12).Value = "Segunda"
$ws.Cells.Item(1146,13).Value = 200
$ws.Cells.Item(1146,14).Value = 2600
$ws.Cells.Item(1146,15).Value = 2700
$ws.Cells.Item(1146,16).Value = 2650
$ws.Cells.Item(1146,17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(1146,18).Value = "Provincia de Limar$([char]0xED)"
$ws.Cells.Item(1146,19).Value = 2650
$ws.Cells.Item(1146,20).Value = 1
